$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 2.1.0 -> 2.2.0-ballot
$meta.Range("B3").Value = "2.2.0-ballot"

# Date: 2025-12-19T08:22:07+00:00 -> 2025-12-19T09:47:21+00:00
$meta.Range("B8").Value = "2025-12-19T09:47:21+00:00"

# Base Definition: append FHIR version pin
$meta.Range("B18").Value = "http://hl7.org/fhir/StructureDefinition/Extension|4.0.1"

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")

# Binding Value Set: append IG version pin
$elements.Range("Z6").Value = "https://interop.esante.gouv.fr/ig/fhir/tddui/ValueSet/tddui-discriminator-vs|2.2.0-ballot"

# Widen column Z (Binding Value Set) to fit the longer value
$elements.Columns.Item(26).ColumnWidth = 67.8
